# "Four new products added" -> a new earnings entry (27/05/2025, RS. 34.72)
# is appended after the existing 26/05/2025 entry, and that existing entry's
# paragraph mark is made bold (matching its bold "13.34" run) in the process.

$d = $word.ActiveDocument

# Locate the paragraph that currently reads "Earnings: RS. 13.34" (the
# 26/05/2025 entry). Range.Text includes the trailing paragraph mark ("`r"),
# so trim that before comparing. We search by text instead of a hard-coded
# paragraph index so the script is resilient to unrelated structural drift.
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    $txt = $p.Range.Text.TrimEnd("`r")
    if ($txt -eq "Earnings: RS. 13.34") {
        $target = $p
        break
    }
}

if ($target -eq $null) {
    throw "Could not find the 'Earnings: RS. 13.34' paragraph"
}

# Replacement content for that paragraph's range:
#  1) the same paragraph, unchanged apart from its paragraph-mark run
#     properties now being bold (b/bCs), matching the bold "13.34" run
#     already present in it;
#  2) a new empty bold paragraph (a spacer, matching the blank-line pattern
#     already used between every other entry in the document);
#  3) a new "Date: 27/05/2025" paragraph;
#  4) a new "Earnings: RS. 34.72" paragraph.
$xmlFrag = '<pkg:xmlData xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage" xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + `
    '<w:p><w:pPr><w:rPr><w:b/><w:bCs/><w:sz w:val="40"/><w:szCs w:val="40"/><w:lang w:val="en-US"/></w:rPr></w:pPr>' + `
        '<w:r><w:rPr><w:sz w:val="40"/><w:szCs w:val="40"/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">Earnings: RS. </w:t></w:r>' + `
        '<w:r><w:rPr><w:b/><w:bCs/><w:sz w:val="40"/><w:szCs w:val="40"/><w:lang w:val="en-US"/></w:rPr><w:t>13.34</w:t></w:r>' + `
    '</w:p>' + `
    '<w:p><w:pPr><w:rPr><w:b/><w:bCs/><w:sz w:val="40"/><w:szCs w:val="40"/><w:lang w:val="en-US"/></w:rPr></w:pPr></w:p>' + `
    '<w:p><w:pPr><w:rPr><w:sz w:val="40"/><w:szCs w:val="40"/><w:lang w:val="en-US"/></w:rPr></w:pPr>' + `
        '<w:r><w:rPr><w:sz w:val="40"/><w:szCs w:val="40"/><w:lang w:val="en-US"/></w:rPr><w:t>Date: 27/05/2025</w:t></w:r>' + `
    '</w:p>' + `
    '<w:p><w:pPr><w:rPr><w:sz w:val="40"/><w:szCs w:val="40"/><w:lang w:val="en-US"/></w:rPr></w:pPr>' + `
        '<w:r><w:rPr><w:sz w:val="40"/><w:szCs w:val="40"/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">Earnings: RS. </w:t></w:r>' + `
        '<w:r><w:rPr><w:b/><w:bCs/><w:sz w:val="40"/><w:szCs w:val="40"/><w:lang w:val="en-US"/></w:rPr><w:t>34.72</w:t></w:r>' + `
    '</w:p>' + `
    '</pkg:xmlData>'

# InsertXML REPLACES the contents of the exact range it's called on, so
# calling it on just this paragraph's Range swaps that one paragraph for the
# four paragraphs above without disturbing anything before or after it.
$target.Range.InsertXML($xmlFrag)

Write-Host "Inserted new earnings entry: Date: 27/05/2025 / Earnings: RS. 34.72"
